# Weekly data refresh: a new daily price record was inserted at row 119
# (Fecha 44474, "Región de Arica y Parinacota", $/atado), pushing all the
# subsequent rows (previously 119-143) down by one to 120-144.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 119; Excel shifts rows 119:143 down to 120:144 and
# carries the date-format style of column D along with the inserted row.
$ws.Rows.Item(119).Insert()

# Populate the newly inserted row 119 with the new market record.
$ws.Cells.Item(119, 1).Value = 10
$ws.Cells.Item(119, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(119, 3).Value = "La Araucanía"
$ws.Cells.Item(119, 4).Value = 44474
$ws.Cells.Item(119, 5).Value = 9
$ws.Cells.Item(119, 6).Value = 100112052
$ws.Cells.Item(119, 7).Value = "Albahaca"
$ws.Cells.Item(119, 8).Value = "Sin especificar"
$ws.Cells.Item(119, 9).Value = "Primera"
$ws.Cells.Item(119, 10).Value = 20
$ws.Cells.Item(119, 11).Value = 7000
$ws.Cells.Item(119, 12).Value = 7000
$ws.Cells.Item(119, 13).Value = 7000
$ws.Cells.Item(119, 14).Value = "$/atado"
$ws.Cells.Item(119, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(119, 16).Value = 7000
$ws.Cells.Item(119, 17).Value = 1
$ws.Cells.Item(119, 18).Value = "Hortaliza"
